$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6449
$wsExhibit.Range("F5").Value = 1300
$wsExhibit.Range("F11").Value = 8030
$wsExhibit.Range("F12").Value = 412
$wsExhibit.Range("F16").Value = 283
$wsExhibit.Range("F20").Value = 305
$wsExhibit.Range("F21").Value = 9982
$wsExhibit.Range("F40").Value = 1871
$wsExhibit.Range("F43").Value = 292
$wsExhibit.Range("F44").Value = 181

# Sheet "全部类型" (index 4 / sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6449
$wsAll.Range("F6").Value = 1300
$wsAll.Range("F13").Value = 8030
$wsAll.Range("F14").Value = 412
$wsAll.Range("F18").Value = 283
$wsAll.Range("F21").Value = 305
$wsAll.Range("F22").Value = 9982
$wsAll.Range("F39").Value = 1872
$wsAll.Range("F43").Value = 292
$wsAll.Range("F44").Value = 181
